$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates.
# NumberFormat is forced to "@" (Text) before writing any value that
# could otherwise be auto-detected as a number, so the cell keeps the
# same text representation as the source feed (e.g. "30.646.37", "1.013").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.646.37"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.124.31"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.61%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.11"
$ws.Range("E5").Value = "  +5.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5277"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4544"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.93"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09091"
$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.181"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.65"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.141.23"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.850"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.093"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.52"
$ws.Range("E16").Value = "  +6.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001176"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06710"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.45"
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.345"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.716.32"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.89"
$ws.Range("E24").Value = "  +3.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.391"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.382.58"
$ws.Range("E26").Value = "  +0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.49"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.27"
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.197"
$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.657"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.365"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.017"
$ws.Range("E35").Value = "  +1.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.159"
$ws.Range("E36").Value = "  +8.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.32"
$ws.Range("E37").Value = "  -1.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02655"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06893"
$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.55"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6927"
$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.279"
$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.80"
$ws.Range("E44").Value = "  +5.00%  "

$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6454"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.772"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000368"
$ws.Range("E48").Value = "  +7.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.257"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.08"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07304"
$ws.Range("E51").Value = "  +2.32%  "

# Rows 28/29: the two coins swapped rank/position - coin name, link,
# price and volume all change together.
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.568"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.77"
$ws.Range("E29").Value = "  +0.68%  "
